$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 18.31647966666667
$ws.Range("H2").Value = 54.949439
$ws.Range("I2").Value = 0.005487334033884006
$ws.Range("J2").Value = 0.005487334033884005
$ws.Range("M2").Value = 2.598166333333333
$ws.Range("N2").Value = 7.794499
$ws.Range("O2").Value = 0.3466013321552429
$ws.Range("P2").Value = 0.3466013321552429
$ws.Range("Q2").Value = 47.58926081511789
$ws.Range("R2").Value = 428.303347336061
$ws.Range("S2").Value = 0.001901917286124999
$ws.Range("T2").Value = 0.001901917286124999
$ws.Range("G3").Value = 18.31647966666667
$ws.Range("H3").Value = 54.949439
$ws.Range("I3").Value = 0.005487334033884006
$ws.Range("J3").Value = 0.005487334033884005
$ws.Range("M3").Value = 4.333403333333333
$ws.Range("O3").Value = 0.5780859172985858
$ws.Range("P3").Value = 0.5780859172985858
$ws.Range("Q3").Value = 79.37269404246555
$ws.Range("R3").Value = 714.3542463821899
$ws.Range("S3").Value = 0.003172150528501585
$ws.Range("T3").Value = 0.003172150528501584
$ws.Range("G4").Value = 18.31647966666667
$ws.Range("H4").Value = 54.949439
$ws.Range("I4").Value = 0.005487334033884006
$ws.Range("J4").Value = 0.005487334033884005
$ws.Range("M4").Value = 0.4692043333333333
$ws.Range("N4").Value = 1.407613
$ws.Range("O4").Value = 0.06259293136852516
$ws.Range("P4").Value = 0.06259293136852516
$ws.Range("Q4").Value = 8.594171631011889
$ws.Range("R4").Value = 77.347544679107
$ws.Range("S4").Value = 0.0003434683225790739
$ws.Range("T4").Value = 0.0003434683225790738
$ws.Range("G5").Value = 18.31647966666667
$ws.Range("H5").Value = 54.949439
$ws.Range("I5").Value = 0.005487334033884006
$ws.Range("J5").Value = 0.005487334033884005
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09534933333333333
$ws.Range("N5").Value = 0.286048
$ws.Range("O5").Value = 0.01271981917764605
$ws.Range("P5").Value = 0.01271981917764604
$ws.Range("Q5").Value = 1.746464125230222
$ws.Range("R5").Value = 15.718177127072
$ws.Range("S5").Value = 0.00006979789667834762
$ws.Range("T5").Value = 0.00006979789667834759
$ws.Range("I6").Value = 0.9472399998689139
$ws.Range("J6").Value = 0.9472399998689137
$ws.Range("M6").Value = 2.598166333333333
$ws.Range("N6").Value = 7.794499
$ws.Range("O6").Value = 0.3466013321552429
$ws.Range("P6").Value = 0.3466013321552429
$ws.Range("Q6").Value = 8215.00042277668
$ws.Range("R6").Value = 73935.00380499012
$ws.Range("S6").Value = 0.3283146458252977
$ws.Range("T6").Value = 0.3283146458252976
$ws.Range("I7").Value = 0.9472399998689139
$ws.Range("J7").Value = 0.9472399998689137
$ws.Range("M7").Value = 4.333403333333333
$ws.Range("O7").Value = 0.5780859172985858
$ws.Range("P7").Value = 0.5780859172985858
$ws.Range("S7").Value = 0.5475861042261334
$ws.Range("T7").Value = 0.5475861042261333
$ws.Range("I8").Value = 0.9472399998689139
$ws.Range("J8").Value = 0.9472399998689137
$ws.Range("M8").Value = 0.4692043333333333
$ws.Range("N8").Value = 1.407613
$ws.Range("O8").Value = 0.06259293136852516
$ws.Range("P8").Value = 0.06259293136852516
$ws.Range("Q8").Value = 1483.551590693122
$ws.Range("R8").Value = 13351.9643162381
$ws.Range("S8").Value = 0.05929052830131671
$ws.Range("T8").Value = 0.0592905283013167
$ws.Range("I9").Value = 0.9472399998689139
$ws.Range("J9").Value = 0.9472399998689137
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.09534933333333333
$ws.Range("N9").Value = 0.286048
$ws.Range("O9").Value = 0.01271981917764605
$ws.Range("P9").Value = 0.01271981917764604
$ws.Range("Q9").Value = 301.4798566186773
$ws.Range("R9").Value = 2713.318709568096
$ws.Range("S9").Value = 0.01204872151616605
$ws.Range("T9").Value = 0.01204872151616605
$ws.Range("G10").Value = 155.6514383333333
$ws.Range("H10").Value = 466.954315
$ws.Range("I10").Value = 0.04663076369111781
$ws.Range("J10").Value = 0.0466307636911178
$ws.Range("M10").Value = 2.598166333333333
$ws.Range("N10").Value = 7.794499
$ws.Range("O10").Value = 0.3466013321552429
$ws.Range("P10").Value = 0.3466013321552429
$ws.Range("Q10").Value = 404.4083268125761
$ws.Range("R10").Value = 3639.674941313185
$ws.Range("S10").Value = 0.01616228481475777
$ws.Range("T10").Value = 0.01616228481475776
$ws.Range("G11").Value = 155.6514383333333
$ws.Range("H11").Value = 466.954315
$ws.Range("I11").Value = 0.04663076369111781
$ws.Range("J11").Value = 0.0466307636911178
$ws.Range("M11").Value = 4.333403333333333
$ws.Range("O11").Value = 0.5780859172985858
$ws.Range("P11").Value = 0.5780859172985858
$ws.Range("Q11").Value = 674.5004617117944
$ws.Range("R11").Value = 6070.50415540615
$ws.Range("S11").Value = 0.02695658780271343
$ws.Range("T11").Value = 0.02695658780271342
$ws.Range("G12").Value = 155.6514383333333
$ws.Range("H12").Value = 466.954315
$ws.Range("I12").Value = 0.04663076369111781
$ws.Range("J12").Value = 0.0466307636911178
$ws.Range("M12").Value = 0.4692043333333333
$ws.Range("N12").Value = 1.407613
$ws.Range("O12").Value = 0.06259293136852516
$ws.Range("P12").Value = 0.06259293136852516
$ws.Range("Q12").Value = 73.03232935556612
$ws.Range("R12").Value = 657.290964200095
$ws.Range("S12").Value = 0.002918756191380052
$ws.Range("T12").Value = 0.002918756191380051
$ws.Range("G13").Value = 155.6514383333333
$ws.Range("H13").Value = 466.954315
$ws.Range("I13").Value = 0.04663076369111781
$ws.Range("J13").Value = 0.0466307636911178
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.09534933333333333
$ws.Range("N13").Value = 0.286048
$ws.Range("O13").Value = 0.01271981917764605
$ws.Range("P13").Value = 0.01271981917764604
$ws.Range("Q13").Value = 14.84126087745778
$ws.Range("R13").Value = 133.57134789712
$ws.Range("S13").Value = 0.0005931348822665613
$ws.Range("T13").Value = 0.0005931348822665611
$ws.Range("G14").Value = 2.142642
$ws.Range("H14").Value = 6.427926
$ws.Range("I14").Value = 0.0006419024060843985
$ws.Range("J14").Value = 0.0006419024060843984
$ws.Range("M14").Value = 2.598166333333333
$ws.Range("N14").Value = 7.794499
$ws.Range("O14").Value = 0.3466013321552429
$ws.Range("P14").Value = 0.3466013321552429
$ws.Range("Q14").Value = 5.566940308786
$ws.Range("R14").Value = 50.102462779074
$ws.Range("S14").Value = 0.0002224842290625082
$ws.Range("T14").Value = 0.0002224842290625082
$ws.Range("G15").Value = 2.142642
$ws.Range("H15").Value = 6.427926
$ws.Range("I15").Value = 0.0006419024060843985
$ws.Range("J15").Value = 0.0006419024060843984
$ws.Range("M15").Value = 4.333403333333333
$ws.Range("O15").Value = 0.5780859172985858
$ws.Range("P15").Value = 0.5780859172985858
$ws.Range("Q15").Value = 9.284931984939998
$ws.Range("R15").Value = 83.56438786446
$ws.Range("S15").Value = 0.0003710747412374689
$ws.Range("T15").Value = 0.0003710747412374688
$ws.Range("G16").Value = 2.142642
$ws.Range("H16").Value = 6.427926
$ws.Range("I16").Value = 0.0006419024060843985
$ws.Range("J16").Value = 0.0006419024060843984
$ws.Range("M16").Value = 0.4692043333333333
$ws.Range("N16").Value = 1.407613
$ws.Range("O16").Value = 0.06259293136852516
$ws.Range("P16").Value = 0.06259293136852516
$ws.Range("Q16").Value = 1.005336911182
$ws.Range("R16").Value = 9.048032200638
$ws.Range("S16").Value = 0.00004017855324933192
$ws.Range("T16").Value = 0.00004017855324933191
$ws.Range("G17").Value = 2.142642
$ws.Range("H17").Value = 6.427926
$ws.Range("I17").Value = 0.0006419024060843985
$ws.Range("J17").Value = 0.0006419024060843984
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.09534933333333333
$ws.Range("N17").Value = 0.286048
$ws.Range("O17").Value = 0.01271981917764605
$ws.Range("P17").Value = 0.01271981917764604
$ws.Range("Q17").Value = 0.204299486272
$ws.Range("R17").Value = 1.838695376448
$ws.Range("S17").Value = 0.000008164882535089472
$ws.Range("T17").Value = 0.000008164882535089469
